# Fix mojibake: cells were double-encoded so the intended "±" (U+00B1)
# character turned into the two-character garble "Â±" (U+00C2 U+00B1).
# This repairs every inline-string cell on the sheets that hold the
# "<value> (<mean> ± <std>)" style result strings: f1_score, training_time
# and test_time. Other sheets (missing_runs, best_seed) are untouched.

$wb = $excel.ActiveWorkbook

$mojibake = [string][char]0x00C2 + [char]0x00B1   # "Â±"
$fixedChar = [string][char]0x00B1                  # "±"

$sheetNames = @("f1_score", "training_time", "test_time")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    for ($r = 1; $r -le 17; $r++) {
        for ($c = 1; $c -le 8; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $text = $cell.Text

            if ($text -like "*$mojibake*") {
                $fixed = $text -replace $mojibake, $fixedChar
                $cell.Value = $fixed
            }
        }
    }
}
